$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " constel·lació, Constel",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "  Constel",
    2
)
